$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing room assignment (column B) and capacity (column E) cells ---
$ws.Range("B6").Value = 'B F2.5'
$ws.Range("B23").Value = 'B F2.2'
$ws.Range("B38").Value = 'A F1.3 - Computer Lab'
$ws.Range("B52").Value = 'B F2.15 - Amphitheater II'
$ws.Range("B60").Value = 'B F2.15 - Amphitheater II'
$ws.Range("B66").Value = 'RC.G4 - GBE IV'
$ws.Range("B83").Value = 'A F1.3 - Computer Lab'
$ws.Range("E83").Value = 25
$ws.Range("B87").Value = 'A B.13 - Class/PSY Lab'
$ws.Range("E87").Value = 12
$ws.Range("B89").Value = 'A B.8 - Fabrication Lab'
$ws.Range("B96").Value = 'A B.2 - EE Lab'
$ws.Range("B97").Value = 'RC1.3 - GSM and Network Laboratories'
$ws.Range("B104").Value = 'RC1.3 - GSM and Network Laboratories'
$ws.Range("B106").Value = 'A B.1 - VACD Multimedia Studio'
$ws.Range("B123").Value = 'B F1.2 - Class/ECON Lab'
$ws.Range("B133").Value = 'A F3.10 - Architecture Classroom'
$ws.Range("B135").Value = 'B F1.35 FBA Conference Room'
$ws.Range("B143").Value = 'B F2.16'
$ws.Range("B153").Value = 'A B.2 - EE Lab'
$ws.Range("B157").Value = 'A F2.16 - Architecture Studio'
$ws.Range("E157").Value = 20
$ws.Range("B162").Value = 'A F2.8 - Drawing Studio & A F2.16 - Architecture Studio'
$ws.Range("B165").Value = 'A B.16 - VACD Drawing Studio'
$ws.Range("B171").Value = 'B F1.35 FBA Conference Room'
$ws.Range("B182").Value = 'A F3.7 - Small Architecture Studio'
$ws.Range("B183").Value = 'B F1.22'
$ws.Range("B187").Value = 'B F2.17'
$ws.Range("B198").Value = 'B F1.23 - Amphitheater I'
$ws.Range("B204").Value = 'RC1.5 - Electronic Laboratory'
$ws.Range("B211").Value = 'A F1.26'
$ws.Range("B218").Value = 'A F2.8 - Drawing Studio & A F2.16 - Architecture Studio'
$ws.Range("B232").Value = 'B F1.24 (MAC Studio)'
$ws.Range("B236").Value = 'A B.2 - EE Lab'
$ws.Range("E236").Value = 25
$ws.Range("B257").Value = 'A B.16 - VACD Drawing Studio'
$ws.Range("B263").Value = 'B F1.10 Class/ART Studio'
$ws.Range("B264").Value = 'B F2.2'
$ws.Range("B280").Value = 'A B.2 - EE Lab'
$ws.Range("B283").Value = 'B F1.2 - Class/ECON Lab'
$ws.Range("B303").Value = 'B F2.27 Creative Writing and Translation Studio'
$ws.Range("E303").Value = 18
$ws.Range("B304").Value = 'RC1.4 - Computer Laboratory'
$ws.Range("E304").Value = 20
$ws.Range("B313").Value = 'A F1.3 - Computer Lab'
$ws.Range("B318").Value = 'RC1.4 - Computer Laboratory'
$ws.Range("E318").Value = 20
$ws.Range("B323").Value = 'A B.13 - Class/PSY Lab'
$ws.Range("E323").Value = 12
$ws.Range("B324").Value = 'A F2.8 - Drawing Studio'
$ws.Range("E324").Value = 25
$ws.Range("B325").Value = 'A B.13 - Class/PSY Lab'
$ws.Range("E325").Value = 12
$ws.Range("B326").Value = 'B F2.27 Creative Writing and Translation Studio'
$ws.Range("E326").Value = 18
$ws.Range("B328").Value = 'B F2.27 Creative Writing and Translation Studio'
$ws.Range("E328").Value = 18
$ws.Range("B331").Value = 'RC1.3 - GSM and Network Laboratories'
$ws.Range("E331").Value = 20
$ws.Range("B333").Value = 'B F2.16'

# --- Append 4 new graduate course rows (334-337) ---
$ws.Range("A334").Value = 'CS600.1'
$ws.Range("B334").Value = 'A B.13 - Class/PSY Lab'
$ws.Range("C334").Value = 'Tue. 17:00 - 19:50'
$ws.Range("D334").Value = 3
$ws.Range("E334").Value = 12
$ws.Range("F334").Value = 'Assigned'

$ws.Range("A335").Value = 'EE603.1'
$ws.Range("B335").Value = 'A B.13 - Class/PSY Lab'
$ws.Range("C335").Value = 'Mon. 15:00 - 17:50'
$ws.Range("D335").Value = 3
$ws.Range("E335").Value = 12
$ws.Range("F335").Value = 'Assigned'

$ws.Range("A336").Value = 'ME605.1'
$ws.Range("B336").Value = 'A B.13 - Class/PSY Lab'
$ws.Range("C336").Value = 'Wed. 17:00 - 19:50'
$ws.Range("D336").Value = 3
$ws.Range("E336").Value = 12
$ws.Range("F336").Value = 'Assigned'

$ws.Range("A337").Value = 'ME580.1'
$ws.Range("B337").Value = 'B F2.27 Creative Writing and Translation Studio'
$ws.Range("C337").Value = 'Mon. 17:00 - 19:50'
$ws.Range("D337").Value = 13
$ws.Range("E337").Value = 18
$ws.Range("F337").Value = 'Assigned'

